# issue #5: property land done
# Renames the Chinese "land" (土地) sheet headers to the canonical English
# field names, fixes a couple of stray full-width spaces/dashes baked into
# some strings, and appends the common metadata columns
# (property_category/category/date/legislator_name/legislator_id/
# source_file/index) that the other sheets already carry.

$wb = $excel.ActiveWorkbook

$wsLand    = $wb.Worksheets.Item(1)   # 土地 (land)
$wsDeposit = $wb.Worksheets.Item(2)   # 存款 (deposit)

# ---------------------------------------------------------------------------
# 1. Sheet "土地" (land) - rename the existing header row to the English
#    field names used across the other sheets.
# ---------------------------------------------------------------------------
$wsLand.Range("B1").Value = "name"
$wsLand.Range("C1").Value = "area"
$wsLand.Range("D1").Value = "share_portion"
$wsLand.Range("E1").Value = "owner"
$wsLand.Range("F1").Value = "register_date"
$wsLand.Range("G1").Value = "register_reason"
$wsLand.Range("H1").Value = "acquire_value"

# ---------------------------------------------------------------------------
# 2. Clean up stray full-width / embedded spaces and dashes in the existing
#    data row.
# ---------------------------------------------------------------------------
$wsLand.Range("B2").Value = "基隆市安樂區大武崙段内寮小段05090010地號"
$wsLand.Range("F2").Value = "91年12月26H"

# ---------------------------------------------------------------------------
# 3. Append the shared metadata columns (I:O) that every other sheet already
#    has. Values are written first, then the header/data formatting is
#    copied over from the existing columns so the new cells pick up the
#    same style indices (bold/centered/bordered header, plain data) instead
#    of minting brand new style entries.
# ---------------------------------------------------------------------------
$wsLand.Range("I1").Value = "property_category"
$wsLand.Range("J1").Value = "category"
$wsLand.Range("K1").Value = "date"
$wsLand.Range("L1").Value = "legislator_name"
$wsLand.Range("M1").Value = "legislator_id"
$wsLand.Range("N1").Value = "source_file"
$wsLand.Range("O1").Value = "index"

$wsLand.Range("I2").Value = "land"
$wsLand.Range("J2").Value = "normal"
$wsLand.Range("K2").Value = "2012-05-01"
$wsLand.Range("L2").Value = "謝國樑"
$wsLand.Range("M2").Value = 1387
$wsLand.Range("N2").Value = "tmpa28e1"
$wsLand.Range("O2").Value = 14

$wsLand.Range("H1").Copy()
$wsLand.Range("I1:O1").PasteSpecial(-4122)

$wsLand.Range("H2").Copy()
$wsLand.Range("I2:O2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Sheet "存款" (deposit) - the same branch-name strings had the same
#    stray embedded spaces; fix them up too.
# ---------------------------------------------------------------------------
$wsDeposit.Range("B2").Value = "基隆市第二信用合作社營業部"
$wsDeposit.Range("B3").Value = "基隆市第二信用合作社營業部"
$wsDeposit.Range("B4").Value = "基隆市第二信用合作社營業部"
$wsDeposit.Range("B5").Value = "基隆市第二信用合作社港東分社"
$wsDeposit.Range("B9").Value = "國泰世華商業銀行板橋分行"
$wsDeposit.Range("B11").Value = "中國信託商業銀行城中分行"
$wsDeposit.Range("B13").Value = "中華郵政股份有限公司基隆"
